# fix melody chord bug
# Rewrite the Duration/Notes pattern pairs (rows 2-21) on the "Pattern"
# worksheet with the corrected values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "0.75_1_1_0.75_0.5"
$ws.Range("B2").Value = "67_71_69_65_67"
$ws.Range("A3").Value = "0.25_0.5_0.25_1_0.25_0.25_1_0.5"
$ws.Range("B3").Value = "69_73_69_69_64_75_64_76"
$ws.Range("A4").Value = "1_0.25_0.25_0.5_2"
$ws.Range("B4").Value = "73_65_68_73_77"
$ws.Range("A5").Value = "0.25_2_0.25_0.5_0.5_0.5"
$ws.Range("B5").Value = "74_74_65_67_67_74"
$ws.Range("A6").Value = "2_0.5_0.75_0.75"
$ws.Range("B6").Value = "71_62_69_60"
$ws.Range("A7").Value = "2_1_0.25_0.75"
$ws.Range("B7").Value = "75_68_73_64"
$ws.Range("A8").Value = "2_0.5_0.75_0.75"
$ws.Range("B8").Value = "77_65_75_65"
$ws.Range("A9").Value = "0.75_2_0.5_0.5_0.25"
$ws.Range("B9").Value = "68_68_67_74_75"
$ws.Range("A10").Value = "2_0.5_0.25_0.25_0.75_0.25"
$ws.Range("B10").Value = "72_60_72_65_71_65"
$ws.Range("A11").Value = "0.25_0.5_1_0.25_2"
$ws.Range("B11").Value = "75_71_75_64_69"
$ws.Range("A12").Value = "1_1_1_0.75_0.25"
$ws.Range("B12").Value = "72_65_72_75_77"
$ws.Range("A13").Value = "0.25_1_0.5_2_0.25"
$ws.Range("B13").Value = "63_72_67_72_72"
$ws.Range("A14").Value = "0.5_0.75_1_1_0.5_0.25"
$ws.Range("B14").Value = "71_60_69_67_60_64"
$ws.Range("A15").Value = "0.5_0.25_1_2_0.25"
$ws.Range("B15").Value = "76_75_64_71_64"
$ws.Range("A16").Value = "2_0.5_1_0.25_0.25"
$ws.Range("B16").Value = "75_73_77_68_65"
$ws.Range("A17").Value = "0.75_0.5_0.5_0.25_2"
$ws.Range("B17").Value = "68_74_67_65_68"
$ws.Range("A18").Value = "0.25_2_1_0.75"
$ws.Range("B18").Value = "69_71_64_60"
$ws.Range("A19").Value = "2_0.5_0.25_0.75_0.5"
$ws.Range("B19").Value = "66_64_69_75_64"
$ws.Range("A20").Value = "2_0.25_0.5_0.25_0.25_0.75"
$ws.Range("B20").Value = "75_65_67_75_75_65"
$ws.Range("A21").Value = "2_1_0.25_0.75"
$ws.Range("B21").Value = "75_67_74_68"
